$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offices = "3.6% CR/LFM+CDL/H:1/Offices`n9.1% CR/LFM+CDL/H:2/Offices`n21.8% CR/LFM+CDL/HBET:3-5/Offices`n1.8% CR/LFM+CDL/HBET:6-/Offices`n3.6% MR/LWAL+CDL/H:1/Offices`n7.3% MR/LWAL+CDL/H:2/Offices`n25.5% MR/LWAL+CDL/HBET:3-5/Offices`n17.3% S/LFM+CDL/HBET:3-5/Offices`n0.9% S/LFM+CDL/HBET:6-/Offices`n8.6% S/LWAL+CDL/HBET:3-5/Offices`n0.5% S/LWAL+CDL/HBET:6-/Offices`n0.0% S/LFM+CDL/H:1/Offices`n0.0% S/LFM+CDL/H:2/Offices`n0.0% S/LWAL+CDL/H:1/Offices`n0.0% S/LWAL+CDL/H:2/Offices`n0.0% W/LFM+CDL/H:1/Offices`n0.0% W/LFM+CDL/H:2/Offices"

$trade = "16.4% CR/LFM+CDL/H:1/Trade`n1.8% CR/LFM+CDL/H:2/Trade`n0.0% CR/LFM+CDL/HBET:3-5/Trade`n0.0% CR/LFM+CDL/HBET:6-/Trade`n32.7% MR/LWAL+CDL/H:1/Trade`n3.6% MR/LWAL+CDL/H:2/Trade`n0.0% MR/LWAL+CDL/HBET:3-5/Trade`n0.0% S/LFM+CDL/HBET:3-5/Trade`n0.0% S/LFM+CDL/HBET:6-/Trade`n0.0% S/LWAL+CDL/HBET:3-5/Trade`n0.0% S/LWAL+CDL/HBET:6-/Trade`n17.8% S/LFM+CDL/H:1/Trade`n0.4% S/LFM+CDL/H:2/Trade`n8.9% S/LWAL+CDL/H:1/Trade`n0.2% S/LWAL+CDL/H:2/Trade`n17.8% W/LFM+CDL/H:1/Trade`n0.4% W/LFM+CDL/H:2/Trade"

$hotels = "4.4% CR/LFM+CDL/H:1/Hotels`n 1.8% CR/LFM+CDL/H:2/Hotels`n 10.6% CR/LFM+CDL/HBET:3-5/Hotels`n 0.9% CR/LFM+CDL/HBET:6-/Hotels`n 3.5% MR/LWAL+CDL/H:1/Hotels`n 7.1% MR/LWAL+CDL/H:2/Hotels`n 25.1% MR/LWAL+CDL/HBET:3-5/Hotels`n 17.2% S/LFM+CDL/HBET:3-5/Hotels`n 2.9% S/LFM+CDL/HBET:6-/Hotels`n 8.4% S/LWAL+CDL/HBET:3-5/Hotels`n 0.4% S/LWAL+CDL/HBET:6-/Hotels`n 0.0% S/LFM+CDL/H:1/Hotels`n 0.0% S/LFM+CDL/H:2/Hotels`n 0.0% S/LWAL+CDL/H:1/Hotels`n 0.0% S/LWAL+CDL/H:2/Hotels`n 14.2% W/LFM+CDL/H:1/Hotels`n 3.5% W/LFM+CDL/H:2/Hotels"

$ws.Range("B2").Value = $offices
$ws.Range("C2").Value = $trade
$ws.Range("D2").Value = $hotels
